$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2 and 3 have their Fecha (D), Volumen (J), Precio minimo (K),
# Precio maximo (L), Precio promedio ponderado (M) and Precio $/Kg (P) swapped.

# Row 2 (new values)
$ws.Range("D2").Value = 44547
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 1600
$ws.Range("M2").Value = 1550
$ws.Range("P2").Value = 1550

# Row 3 (new values)
$ws.Range("D3").Value = 44875
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 1600
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = 1650
$ws.Range("P3").Value = 1650
